# Update the "NEW" map data sheet:
#  - remove the resolved case (old row 22: Caso 6173 "Picada volvio a entrar como caso 6325")
#  - split the old "Reparar rienda" / "Picada" entry pair into two distinct new cases
#    (Caso 6229 @ ALVAREZ THOMAS AV. 309, and Caso 6228 @ NEWBERY, JORGE AV. 3416),
#    inserted where old row 33 ("Santa maria de oro 2722") used to begin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the now-resolved row (shifts everything below it up by one).
$ws.Rows(22).Delete()

# 2) Insert two fresh blank rows at 32:33 (shifts everything from the old
#    row 32 onward back down by two).
$ws.Rows("32:33").Insert()

# 3) Populate the two new rows. Columns A, B, D, E hold values that look like
#    numbers/dates, so force a text number format first so they are stored
#    the same way as the rest of the sheet (plain text), not auto-converted.
$textCols = 1,2,4,5

foreach ($col in $textCols) {
    $ws.Cells.Item(32, $col).NumberFormat = "@"
    $ws.Cells.Item(33, $col).NumberFormat = "@"
}

# Row 32: Caso 6229
$ws.Cells.Item(32,1).Value2  = "6229"
$ws.Cells.Item(32,2).Value2  = "6/24/2025"
$ws.Cells.Item(32,3).Value2  = "ALVAREZ THOMAS AV. 309"
$ws.Cells.Item(32,4).Value2  = "13"
$ws.Cells.Item(32,5).Value2  = "807762987"
$ws.Cells.Item(32,6).Value2  = "NEW"
$ws.Cells.Item(32,7).Value2  = "Pendiente"
$ws.Cells.Item(32,8).Value2  = "Reparar rienda "
$ws.Cells.Item(32,9).Value2  = 1
$ws.Cells.Item(32,10).Value2 = "Tensor"
$ws.Cells.Item(32,11).Value2 = "Sin equipos"
$ws.Cells.Item(32,12).Value2 = "Terminal"
$ws.Cells.Item(32,13).Value2 = -58.44848
$ws.Cells.Item(32,14).Value2 = -34.581338
$ws.Cells.Item(32,15).Value2 = "Palermo"
$ws.Cells.Item(32,16).Value2 = "Capital Sur"

# Row 33: Caso 6228
$ws.Cells.Item(33,1).Value2  = "6228"
$ws.Cells.Item(33,2).Value2  = "6/24/2025"
$ws.Cells.Item(33,3).Value2  = "NEWBERY, JORGE AV. 3416"
$ws.Cells.Item(33,4).Value2  = "13"
$ws.Cells.Item(33,5).Value2  = "807762990"
$ws.Cells.Item(33,6).Value2  = "NEW"
$ws.Cells.Item(33,7).Value2  = "Pendiente"
$ws.Cells.Item(33,8).Value2  = "Reparar rienda"
$ws.Cells.Item(33,9).Value2  = 1
$ws.Cells.Item(33,10).Value2 = "Tensor"
$ws.Cells.Item(33,11).Value2 = "Sin equipos"
$ws.Cells.Item(33,12).Value2 = "Terminal"
$ws.Cells.Item(33,13).Value2 = -58.448496
$ws.Cells.Item(33,14).Value2 = -34.58182
$ws.Cells.Item(33,15).Value2 = "Colegiales"
$ws.Cells.Item(33,16).Value2 = "Capital Norte"
